$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.892.42'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '3.938.85'
$ws.Range("E3").Value = '  +7.11%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.95'
$ws.Range("E5").Value = '  -0.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.71'
$ws.Range("E6").Value = '  -1.45%  '

$ws.Range("D7").Value = '3.937.93'
$ws.Range("E7").Value = '  +7.18%  '

$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("E10").Value = '  +2.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.42'
$ws.Range("E11").Value = '  +2.61%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.486'
$ws.Range("E12").Value = '  +1.72%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.11'
$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("E14").Value = '  +1.49%  '

$ws.Range("D15").Value = '4.588.45'
$ws.Range("E15").Value = '  +6.22%  '

$ws.Range("D16").Value = '3.933.23'
$ws.Range("E16").Value = '  +6.52%  '

$ws.Range("D17").Value = '69.963.00'
$ws.Range("E17").Value = '  +0.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.51'
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("E19").Value = '  -2.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.71'
$ws.Range("E20").Value = '  +2.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '508.33'
$ws.Range("E21").Value = '  +2.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.84'
$ws.Range("E22").Value = '  +7.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.757'
$ws.Range("E23").Value = '  +6.86%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.16'
$ws.Range("E24").Value = '  +1.98%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.45'
$ws.Range("E25").Value = '  -4.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000140'
$ws.Range("E26").Value = '  +6.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.63'
$ws.Range("E27").Value = '  -2.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.42'
$ws.Range("E28").Value = '  -7.84%  '

$ws.Range("E29").Value = '  +0.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.58'
$ws.Range("E30").Value = '  +5.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.01'
$ws.Range("E31").Value = '  +3.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '34.10'
$ws.Range("E32").Value = '  +13.30%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.90'
$ws.Range("E33").Value = '  +1.72%  '

$ws.Range("E34").Value = '  +0.88%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.05'
$ws.Range("E36").Value = '  +1.28%  '

$ws.Range("E37").Value = '  +2.50%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.142'
$ws.Range("E38").Value = '  +3.42%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '475.13'
$ws.Range("E39").Value = '  +12.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.335'
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.05'
$ws.Range("E41").Value = '  -0.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '49.81'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.93'
$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("E44").Value = '  +0.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.76'
$ws.Range("E45").Value = '  -5.22%  '

$ws.Range("D46").Value = '2.963.07'
$ws.Range("E46").Value = '  +0.74%  '

$ws.Range("E47").Value = '  +2.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.54'
$ws.Range("E48").Value = '  +1.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.77'
$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '24.98'
$ws.Range("E51").Value = '  +18.09%  '
